$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the "Native" column (G) for several species rows where it had
# previously been recorded as "Unknown" / "Native/Unknown".
$ws.Range("G5").Value2  = "Native"
$ws.Range("G29").Value2 = "Native"
$ws.Range("G45").Value2 = "Native"
$ws.Range("G46").Value2 = "Native"
$ws.Range("G61").Value2 = "Native"
$ws.Range("G62").Value2 = "Native"
$ws.Range("G63").Value2 = "Native"
$ws.Range("G70").Value2 = "Native"
$ws.Range("G83").Value2 = "Native"

# Rename the UNGR1 entry (row 82) to clarify it is likely BOER4 or Aristida.
$ws.Range("D82").Value2 = "UNGR1 (BOER/Aristida)"
$ws.Range("E82").Value2 = "UNGR1 (BOER/Aristida).SRER"
$ws.Range("F82").Value2 = "Unknown grass 1, BOER4 or Aristida, SRER"

# Reflect where the user's selection/scroll position ended up after editing.
$ws.Range("A71").Select()
$ws.Range("D82:F82").Select()
